$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 values to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 6.38
$ws.Range("C5").Value = 4.89
$ws.Range("D5").Value = 0.39
$ws.Range("E5").Value = 13.41
$ws.Range("F5").Value = 11.48
$ws.Range("G5").Value = 4.77
$ws.Range("I5").Value = 7.31
$ws.Range("J5").Value = 3.42
$ws.Range("K5").Value = 5.24
$ws.Range("L5").Value = 5.53
$ws.Range("M5").Value = 5.7
$ws.Range("N5").Value = 1.55
$ws.Range("O5").Value = 4.77
$ws.Range("P5").Value = 7.02
$ws.Range("Q5").Value = 4.02
$ws.Range("R5").Value = 0.36
$ws.Range("S5").Value = 0.06
$ws.Range("T5").Value = 66.94
$ws.Range("U5").Value = 13.72
$ws.Range("V5").Value = 4.36
$ws.Range("W5").Value = 9.289999999999999
$ws.Range("X5").Value = 5.08
$ws.Range("Y5").Value = 0.65
$ws.Range("AA5").Value = 3.66
$ws.Range("AB5").Value = 4
$ws.Range("AC5").Value = 3.94
$ws.Range("AD5").Value = 5.82
$ws.Range("AF5").Value = 22.59
$ws.Range("AG5").Value = 2.72
$ws.Range("AH5").Value = 5.48

# Remove row 6 entirely (data trimmed to 1000 rows / fewer samples)
$ws.Rows(6).Delete()
